$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 269.54544
$ws.Range("I18").Value = 269.54544
$ws.Range("K18").Value = 269.54544
$ws.Range("M18").Value = 14.45456000000001
$ws.Range("H33").Value = 407.9091
$ws.Range("I33").Value = 407.9091
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 407.9091
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H86").Value = 4980.5
$ws.Range("I86").Value = 4980
$ws.Range("J86").Value = 4981
$ws.Range("K86").Value = 4980
$ws.Range("L86").Value = 4981
$ws.Range("M86").Value = -3857
$ws.Range("N86").Value = -7227
$ws.Range("H89").Value = 4980.5
$ws.Range("I89").Value = 4980
$ws.Range("J89").Value = 4981
$ws.Range("K89").Value = 24930
$ws.Range("L89").Value = 24905
$ws.Range("M89").Value = -19284
$ws.Range("N89").Value = -36137
$ws.Range("H98").Value = 603.7037
$ws.Range("I98").Value = 768.875
$ws.Range("J98").Value = 363.45456
$ws.Range("K98").Value = 768.875
$ws.Range("L98").Value = 363.45456
$ws.Range("M98").Value = 729.125
$ws.Range("N98").Value = -3359.45456
$ws.Range("H122").Value = 603.7037
$ws.Range("I122").Value = 768.875
$ws.Range("J122").Value = 363.45456
$ws.Range("K122").Value = 2306.625
$ws.Range("L122").Value = 1090.36368
$ws.Range("M122").Value = 143.375
$ws.Range("N122").Value = -5990.36368
$ws.Range("H127").Value = 1999.5
$ws.Range("J127").Value = 1999
$ws.Range("L127").Value = 5997
$ws.Range("N127").Value = -15917
$ws.Range("H129").Value = 2629.2222
$ws.Range("I129").Value = 2338.2
$ws.Range("J129").Value = 2993
$ws.Range("K129").Value = 7014.599999999999
$ws.Range("L129").Value = 8979
$ws.Range("M129").Value = -2014.599999999999
$ws.Range("N129").Value = -18979
$ws.Range("H131").Value = 8559.799999999999
$ws.Range("I131").Value = 8449.75
$ws.Range("K131").Value = 25349.25
$ws.Range("M131").Value = -20309.25
$ws.Range("H138").Value = 3961.5557
$ws.Range("I138").Value = 3814.2307
$ws.Range("J138").Value = 4044.8262
$ws.Range("K138").Value = 11442.6921
$ws.Range("L138").Value = 12134.4786
$ws.Range("M138").Value = -6302.6921
$ws.Range("N138").Value = -22414.4786
$ws.Range("H141").Value = 6271.143
$ws.Range("J141").Value = 5579.8
$ws.Range("L141").Value = 16739.4
$ws.Range("N141").Value = -27099.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7276.1313
$ws.Range("I32").Value = 6916.6807
$ws.Range("K32").Value = 6916.6807
$ws.Range("M32").Value = -6629.6807
$ws.Range("H61").Value = 1362.5294
$ws.Range("I61").Value = 962.03705
$ws.Range("K61").Value = 962.03705
$ws.Range("M61").Value = -750.03705
$ws.Range("H97").Value = 499.8889
$ws.Range("I97").Value = 315.66666
$ws.Range("J97").Value = 868.3333
$ws.Range("K97").Value = 315.66666
$ws.Range("L97").Value = 868.3333
$ws.Range("M97").Value = 180.33334
$ws.Range("N97").Value = -1860.3333
$ws.Range("H122").Value = 25000
$ws.Range("I122").Value = 25000
$ws.Range("K122").Value = 75000
$ws.Range("M122").Value = -72550
$ws.Range("H132").Value = 2913.0908
$ws.Range("I132").Value = 2149.4285
$ws.Range("J132").Value = 4249.5
$ws.Range("K132").Value = 6448.2855
$ws.Range("L132").Value = 12748.5
$ws.Range("M132").Value = -3918.2855
$ws.Range("N132").Value = -17808.5
$ws.Range("H136").Value = 1362.5294
$ws.Range("I136").Value = 962.03705
$ws.Range("K136").Value = 2886.11115
$ws.Range("M136").Value = -336.1111500000002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 833
$ws.Range("I22").Value = 833
$ws.Range("K22").Value = 833
$ws.Range("M22").Value = -660
$ws.Range("H105").Value = 4854.7144
$ws.Range("I105").Value = 4663.8335
$ws.Range("J105").Value = 6000
$ws.Range("K105").Value = 4663.8335
$ws.Range("L105").Value = 6000
$ws.Range("M105").Value = -2916.8335
$ws.Range("N105").Value = -9494
$ws.Range("H134").Value = 2553.84
$ws.Range("I134").Value = 2220.2727
$ws.Range("K134").Value = 6660.8181
$ws.Range("M134").Value = -4125.8181

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 405.5
$ws.Range("I7").Value = 49.4
$ws.Range("K7").Value = 49.4
$ws.Range("M7").Value = 63.6
$ws.Range("H16").Value = 1207.1428
$ws.Range("I16").Value = 1162.5
$ws.Range("J16").Value = 1266.6666
$ws.Range("K16").Value = 1162.5
$ws.Range("L16").Value = 1266.6666
$ws.Range("M16").Value = -875.5
$ws.Range("N16").Value = -1840.6666
$ws.Range("H28").Value = 35000
$ws.Range("J28").Value = 35000
$ws.Range("L28").Value = 35000
$ws.Range("N28").Value = -35490
$ws.Range("H58").Value = 2929.6428
$ws.Range("I58").Value = 2871.75
$ws.Range("K58").Value = 2871.75
$ws.Range("M58").Value = -2668.75
$ws.Range("H113").Value = 1207.1428
$ws.Range("I113").Value = 1162.5
$ws.Range("J113").Value = 1266.6666
$ws.Range("K113").Value = 1162.5
$ws.Range("L113").Value = 1266.6666
$ws.Range("M113").Value = 1007.5
$ws.Range("N113").Value = -5606.6666
$ws.Range("H136").Value = 2929.6428
$ws.Range("I136").Value = 2871.75
$ws.Range("K136").Value = 8615.25
$ws.Range("M136").Value = -6065.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 283
$ws.Range("J23").Value = 283
$ws.Range("L23").Value = 849
$ws.Range("N23").Value = -1319
$ws.Range("H69").Value = 1567.375
$ws.Range("I69").Value = 1307.8
$ws.Range("J69").Value = 2000
$ws.Range("K69").Value = 3923.4
$ws.Range("L69").Value = 6000
$ws.Range("M69").Value = -3112.4
$ws.Range("N69").Value = -7622
$ws.Range("H72").Value = 1567.375
$ws.Range("I72").Value = 1307.8
$ws.Range("J72").Value = 2000
$ws.Range("K72").Value = 11770.2
$ws.Range("L72").Value = 18000
$ws.Range("M72").Value = -7714.199999999999
$ws.Range("N72").Value = -26112
$ws.Range("H107").Value = 1668.4286
$ws.Range("J107").Value = 279.83334
$ws.Range("L107").Value = 839.5000200000001
$ws.Range("N107").Value = -4679.50002
$ws.Range("H112").Value = 1349.5
$ws.Range("I112").Value = 900
$ws.Range("J112").Value = 1799
$ws.Range("K112").Value = 2700
$ws.Range("L112").Value = 5397
$ws.Range("M112").Value = -1592
$ws.Range("N112").Value = -7613
$ws.Range("H121").Value = 2095.2
$ws.Range("J121").Value = 2479
$ws.Range("L121").Value = 7437
$ws.Range("N121").Value = -10057
$ws.Range("H122").Value = 112833
$ws.Range("J122").Value = 112833
$ws.Range("L122").Value = 1015497
$ws.Range("N122").Value = -1020397
$ws.Range("H132").Value = 3770.6667
$ws.Range("J132").Value = 3548.0715
$ws.Range("L132").Value = 31932.6435
$ws.Range("N132").Value = -36992.6435

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5112.4287
$ws.Range("I80").Value = 3946.75
$ws.Range("K80").Value = 3946.75
$ws.Range("M80").Value = -2948.75
$ws.Range("H83").Value = 5112.4287
$ws.Range("I83").Value = 3946.75
$ws.Range("K83").Value = 19733.75
$ws.Range("M83").Value = -14741.75
$ws.Range("H97").Value = 674.7273
$ws.Range("I97").Value = 602
$ws.Range("K97").Value = 602
$ws.Range("M97").Value = -106
$ws.Range("H107").Value = 584
$ws.Range("I107").Value = 557
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 557
$ws.Range("L107").Value = 800
$ws.Range("M107").Value = 1363
$ws.Range("N107").Value = -4640
$ws.Range("H136").Value = 74992.664
$ws.Range("J136").Value = 74992.664
$ws.Range("L136").Value = 224977.992
$ws.Range("N136").Value = -230077.992

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 400.5
$ws.Range("I55").Value = 404.6
$ws.Range("K55").Value = 404.6
$ws.Range("M55").Value = -231.6
$ws.Range("H82").Value = 1634.7142
$ws.Range("I82").Value = 1496.2
$ws.Range("K82").Value = 1496.2
$ws.Range("M82").Value = -1135.2
$ws.Range("H85").Value = 1634.7142
$ws.Range("I85").Value = 1496.2
$ws.Range("K85").Value = 1496.2
$ws.Range("M85").Value = -248.2
$ws.Range("H130").Value = 50555.332
$ws.Range("J130").Value = 50555.332
$ws.Range("L130").Value = 50555.332
$ws.Range("N130").Value = -60595.332
$ws.Range("H136").Value = 3306.9714
$ws.Range("I136").Value = 3636.8
$ws.Range("J136").Value = 2482.4
$ws.Range("K136").Value = 10910.4
$ws.Range("L136").Value = 7447.200000000001
$ws.Range("M136").Value = -8360.400000000001
$ws.Range("N136").Value = -12547.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4750
$ws.Range("I122").Value = 4750
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 14250
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("N133").ClearContents()
